# msz - first smoke test is running
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: add the "id=..." control identifiers next to the existing A1 marker
# ---------------------------------------------------------------------------
$row1 = @{
    "B1" = "id=make"
    "C1" = "id=engineperformance"
    "D1" = "id=dateofmanufacture"
    "E1" = "id=numberofseats"
    "F1" = "id=fuel"
    "I1" = "id=listprice"
    "J1" = "id=licenseplatenumber"
    "K1" = "id=annualmileage"
    "L1" = "id=nextenterinsurantdata"
}
foreach ($addr in $row1.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $row1[$addr]
}
# G1 / H1 keep the default (General) style
$ws.Range("G1").Value = "id=payload"
$ws.Range("H1").Value = "id=totalweight"

# ---------------------------------------------------------------------------
# Row 2: highlighted (yellow) but otherwise empty helper row
# ---------------------------------------------------------------------------
$ws.Range("B2:F2").NumberFormat = "@"
$ws.Range("B2:F2").Interior.Color = 65535
$ws.Range("I2:L2").NumberFormat = "@"
$ws.Range("I2:L2").Interior.Color = 65535
$ws.Range("G2:H2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Row 3: rename the combo-box control references to the new selector names,
# and make the payload / total weight headers match the rest of the row
# (same grey header fill/bold as the other header cells, not the yellow one)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "selMake"
$ws.Range("E3").Value = "selNumberOfSeats"
$ws.Range("F3").Value = "selFuelType"
$ws.Range("A3").Copy()
$ws.Range("G3:H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 5: new smoke-test data row
# ---------------------------------------------------------------------------
$ws.Range("A5:B5").NumberFormat = "@"
$ws.Range("A5").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("B5").Value = "Audi"

$ws.Range("C5").Value = 100
$ws.Range("A1").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("D5:L5").NumberFormat = "@"
$ws.Range("D5").Value = "11/29/2011"
$ws.Range("E5").Value = "5"
$ws.Range("F5").Value = "Petrol"
$ws.Range("G5").Value = "1000"
$ws.Range("H5").Value = "1000"
$ws.Range("I5").Value = "25000"
$ws.Range("J5").Value = "12345ABCDE"
$ws.Range("K5").Value = "12500"
$ws.Range("L5").Value = "X"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 34.1
$ws.Columns("L").ColumnWidth = 21

# ---------------------------------------------------------------------------
# Selection (active cell) left at T19 by the author
# ---------------------------------------------------------------------------
$ws.Range("T19").Select()

$wb.Save()
